# Auto-generated edit script: updates market-price derived columns
# (currentAveragePrice/HQ, LevePriceHQ, LeveProfitHQ, etc.) across several sheets,
# per the scheduled runner's refreshed pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 44250
$ws.Range("J75").Value = 44250
$ws.Range("L75").Value = 44250
$ws.Range("N75").Value = -46122

$ws.Range("H78").Value = 44250
$ws.Range("J78").Value = 44250
$ws.Range("L78").Value = 132750
$ws.Range("N78").Value = -142110

$ws.Range("H87").Value = 16744
$ws.Range("J87").Value = 16744
$ws.Range("L87").Value = 16744
$ws.Range("N87").Value = -19240

$ws.Range("H90").Value = 16744
$ws.Range("J90").Value = 16744
$ws.Range("L90").Value = 50232
$ws.Range("N90").Value = -62712

$ws.Range("H114").Value = 42805.332
$ws.Range("J114").Value = 42805.332
$ws.Range("L114").Value = 42805.332
$ws.Range("N114").Value = -51483.332

$ws.Range("H123").Value = 37483.5
$ws.Range("J123").Value = 37483.5
$ws.Range("L123").Value = 37483.5
$ws.Range("N123").Value = -47283.5

$ws.Range("H133").Value = 56859.6
$ws.Range("J133").Value = 56859.6
$ws.Range("L133").Value = 56859.6
$ws.Range("N133").Value = -66979.60000000001

$ws.Range("H137").Value = 4343.275
$ws.Range("I137").Value = 1086.3636
$ws.Range("J137").Value = 8323.944
$ws.Range("K137").Value = 3259.0908
$ws.Range("L137").Value = 24971.832
$ws.Range("M137").Value = -709.0907999999999
$ws.Range("N137").Value = -30071.832

$ws.Range("H138").Value = 1525.9697
$ws.Range("I138").Value = 751.3570999999999
$ws.Range("J138").Value = 2096.7368
$ws.Range("K138").Value = 2254.0713
$ws.Range("L138").Value = 6290.2104
$ws.Range("M138").Value = 2885.9287
$ws.Range("N138").Value = -16570.2104

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 38440.223
$ws.Range("J80").Value = 38440.223
$ws.Range("L80").Value = 38440.223
$ws.Range("N80").Value = -40436.223

$ws.Range("H83").Value = 38440.223
$ws.Range("J83").Value = 38440.223
$ws.Range("L83").Value = 115320.669
$ws.Range("N83").Value = -125304.669

$ws.Range("H113").Value = 46992
$ws.Range("J113").Value = 46992
$ws.Range("L113").Value = 46992
$ws.Range("N113").Value = -55670

$ws.Range("H128").Value = 50421
$ws.Range("J128").Value = 50421
$ws.Range("L128").Value = 50421
$ws.Range("N128").Value = -60381

$ws.Range("H131").Value = 44379
$ws.Range("J131").Value = 44379
$ws.Range("L131").Value = 44379
$ws.Range("N131").Value = -54459

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 26109.445
$ws.Range("J132").Value = 26109.445
$ws.Range("L132").Value = 26109.445
$ws.Range("N132").Value = -36229.445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 47924
$ws.Range("J20").Value = 47924
$ws.Range("L20").Value = 47924
$ws.Range("N20").Value = -48396

$ws.Range("H30").Value = 47924
$ws.Range("J30").Value = 47924
$ws.Range("L30").Value = 47924
$ws.Range("N30").Value = -48106

$ws.Range("H31").Value = 2723.87
$ws.Range("I31").Value = 879.6429000000001
$ws.Range("J31").Value = 3441.0693
$ws.Range("K31").Value = 879.6429000000001
$ws.Range("L31").Value = 3441.0693
$ws.Range("M31").Value = -584.6429000000001
$ws.Range("N31").Value = -4031.0693

$ws.Range("H34").Value = 2723.87
$ws.Range("I34").Value = 879.6429000000001
$ws.Range("J34").Value = 3441.0693
$ws.Range("K34").Value = 879.6429000000001
$ws.Range("L34").Value = 3441.0693
$ws.Range("M34").Value = -677.6429000000001
$ws.Range("N34").Value = -3845.0693

$ws.Range("H110").Value = 46348.5
$ws.Range("J110").Value = 46348.5
$ws.Range("L110").Value = 46348.5
$ws.Range("N110").Value = -54528.5

$ws.Range("H128").Value = 47924
$ws.Range("J128").Value = 47924
$ws.Range("L128").Value = 47924
$ws.Range("N128").Value = -57884

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

$ws.Range("H131").Value = 41892
$ws.Range("J131").Value = 41892
$ws.Range("L131").Value = 41892
$ws.Range("N131").Value = -51972

$ws.Range("H137").Value = 52499.92
$ws.Range("J137").Value = 52499.92
$ws.Range("L137").Value = 52499.92
$ws.Range("N137").Value = -62699.92

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 50855
$ws.Range("J114").Value = 50855
$ws.Range("L114").Value = 50855
$ws.Range("N114").Value = -59533

$ws.Range("H116").Value = 48742
$ws.Range("J116").Value = 48742
$ws.Range("L116").Value = 48742
$ws.Range("N116").Value = -57920

$ws.Range("H124").Value = 42776
$ws.Range("J124").Value = 42776
$ws.Range("L124").Value = 42776
$ws.Range("N124").Value = -52596

$ws.Range("H130").Value = 53992
$ws.Range("J130").Value = 53992
$ws.Range("L130").Value = 53992
$ws.Range("N130").Value = -64032

$ws.Range("H132").Value = 2337.9348
$ws.Range("I132").Value = 1532.5172
$ws.Range("J132").Value = 3711.8823
$ws.Range("K132").Value = 4597.5516
$ws.Range("L132").Value = 11135.6469
$ws.Range("M132").Value = -2067.5516
$ws.Range("N132").Value = -16195.6469

$ws.Range("H135").Value = 47925
$ws.Range("J135").Value = 47925
$ws.Range("L135").Value = 47925
$ws.Range("N135").Value = -58065

$ws.Range("H139").Value = 72494
$ws.Range("J139").Value = 72494
$ws.Range("L139").Value = 72494
$ws.Range("N139").Value = -82774

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 43175.668
$ws.Range("J88").Value = 43175.668
$ws.Range("L88").Value = 43175.668
$ws.Range("N88").Value = -44031.668

$ws.Range("H91").Value = 43175.668
$ws.Range("J91").Value = 43175.668
$ws.Range("L91").Value = 43175.668
$ws.Range("N91").Value = -46139.668

$ws.Range("H96").Value = 37598.5
$ws.Range("J96").Value = 37598.5
$ws.Range("L96").Value = 37598.5
$ws.Range("N96").Value = -43090.5

$ws.Range("H99").Value = 21958.889
$ws.Range("I99").Value = 11526
$ws.Range("K99").Value = 11526
$ws.Range("M99").Value = -8531

$ws.Range("H102").Value = 49561
$ws.Range("J102").Value = 49561
$ws.Range("L102").Value = 49561
$ws.Range("N102").Value = -56051

$ws.Range("H116").Value = 49676
$ws.Range("J116").Value = 49676
$ws.Range("L116").Value = 49676
$ws.Range("N116").Value = -58854

$ws.Range("H123").Value = 39425
$ws.Range("J123").Value = 39425
$ws.Range("L123").Value = 39425
$ws.Range("N123").Value = -49225

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = $null

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null

$ws.Range("H128").Value = 48429
$ws.Range("J128").Value = 48429
$ws.Range("L128").Value = 48429
$ws.Range("N128").Value = -58389

$ws.Range("H137").Value = 36548
$ws.Range("J137").Value = 36548
$ws.Range("L137").Value = 36548
$ws.Range("N137").Value = -46748

$ws.Range("H139").Value = 48933
$ws.Range("J139").Value = 48933
$ws.Range("L139").Value = 48933
$ws.Range("N139").Value = -59213

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value = 29270.285
$ws.Range("J106").Value = 29270.285
$ws.Range("L106").Value = 29270.285
$ws.Range("N106").Value = -31794.285

$ws.Range("H131").Value = 50178.668
$ws.Range("J131").Value = 50178.668
$ws.Range("L131").Value = 50178.668
$ws.Range("N131").Value = -60258.668

$ws.Range("H139").Value = 55700
$ws.Range("J139").Value = 55700
$ws.Range("L139").Value = 55700
$ws.Range("N139").Value = -65980
